$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Exempt under provisions of 33 ILCS 200/31-45, Paragraph ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Exempt under provisions of 35 ILCS 200/31-45, Paragraph ", 2
)
